$d = $word.ActiveDocument

$d.Content.Find.Execute("2026-01-28 Wednesday", $true, $false, $false, $false, $false, $true, 1, $false, "2026-01-29 Thursday", 2) | Out-Null
$d.Content.Find.Execute("951×2=1902", $true, $false, $false, $false, $false, $true, 1, $false, "512×3=1536", 2) | Out-Null
$d.Content.Find.Execute("817×5=4085", $true, $false, $false, $false, $false, $true, 1, $false, "221×2=442", 2) | Out-Null
$d.Content.Find.Execute("146×2=292", $true, $false, $false, $false, $false, $true, 1, $false, "206×4=824", 2) | Out-Null
$d.Content.Find.Execute("938×2=1876", $true, $false, $false, $false, $false, $true, 1, $false, "543×4=2172", 2) | Out-Null
$d.Content.Find.Execute("190×2=380", $true, $false, $false, $false, $false, $true, 1, $false, "775×8=6200", 2) | Out-Null
$d.Content.Find.Execute("244×6=1464", $true, $false, $false, $false, $false, $true, 1, $false, "991×5=4955", 2) | Out-Null
$d.Content.Find.Execute("182×7=1274", $true, $false, $false, $false, $false, $true, 1, $false, "730×6=4380", 2) | Out-Null
$d.Content.Find.Execute("861×8=6888", $true, $false, $false, $false, $false, $true, 1, $false, "655×4=2620", 2) | Out-Null
$d.Content.Find.Execute("203×4=812", $true, $false, $false, $false, $false, $true, 1, $false, "333×7=2331", 2) | Out-Null
$d.Content.Find.Execute("103×6=618", $true, $false, $false, $false, $false, $true, 1, $false, "648×8=5184", 2) | Out-Null
$d.Content.Find.Execute("509×4=2036", $true, $false, $false, $false, $false, $true, 1, $false, "257×7=1799", 2) | Out-Null
$d.Content.Find.Execute("827×4=3308", $true, $false, $false, $false, $false, $true, 1, $false, "497×8=3976", 2) | Out-Null
$d.Content.Find.Execute("919×9=8271", $true, $false, $false, $false, $false, $true, 1, $false, "604×7=4228", 2) | Out-Null
$d.Content.Find.Execute("273×2=546", $true, $false, $false, $false, $false, $true, 1, $false, "739×3=2217", 2) | Out-Null
$d.Content.Find.Execute("235×7=1645", $true, $false, $false, $false, $false, $true, 1, $false, "436×4=1744", 2) | Out-Null
$d.Content.Find.Execute("693×2=1386", $true, $false, $false, $false, $false, $true, 1, $false, "451×8=3608", 2) | Out-Null
$d.Content.Find.Execute("604×8=4832", $true, $false, $false, $false, $false, $true, 1, $false, "986×7=6902", 2) | Out-Null
$d.Content.Find.Execute("684×2=1368", $true, $false, $false, $false, $false, $true, 1, $false, "659×8=5272", 2) | Out-Null
$d.Content.Find.Execute("376×4=1504", $true, $false, $false, $false, $false, $true, 1, $false, "530×6=3180", 2) | Out-Null
$d.Content.Find.Execute("618×9=5562", $true, $false, $false, $false, $false, $true, 1, $false, "738×3=2214", 2) | Out-Null
$d.Content.Find.Execute("753×7=5271", $true, $false, $false, $false, $false, $true, 1, $false, "346×8=2768", 2) | Out-Null
$d.Content.Find.Execute("620×7=4340", $true, $false, $false, $false, $false, $true, 1, $false, "988×6=5928", 2) | Out-Null
$d.Content.Find.Execute("300×3=900", $true, $false, $false, $false, $false, $true, 1, $false, "721×3=2163", 2) | Out-Null
$d.Content.Find.Execute("960×7=6720", $true, $false, $false, $false, $false, $true, 1, $false, "608×6=3648", 2) | Out-Null
$d.Content.Find.Execute("626×6=3756", $true, $false, $false, $false, $false, $true, 1, $false, "254×6=1524", 2) | Out-Null
